$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A (shifts RollNo/Name/Class/Address to B:E)
$ws.Columns("A:A").Insert()

# New column header + values
$ws.Range("A1").Value = "TC_ID"
$ws.Range("A2").Value = "TC_1"
$ws.Range("A3").Value = "TC_2"
$ws.Range("A4").Value = "TC_3"
$ws.Range("A5").Value = "TC_4"
$ws.Range("A6").Value = "TC_5"

$ws.Range("A3:A6").Select() | Out-Null
